$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.913.01"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "1.876.32"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  -0.88%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4842"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3813"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07370"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9411"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  +3.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07783"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "1.889.51"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.546"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.599"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008871"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").Value = "27.905.01"
$ws.Range("E20").Value = "  +2.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "2.118.96"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.52%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.932"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.041"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.972"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08893"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.227"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7724"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.647"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.726"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.128"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02048"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("E39").Value = "  +5.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05378"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.22%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.048"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.537"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1528"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4883"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.668"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06116"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.57%  "
